$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.766.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.84%  '
$ws.Range("D3").Value = "'2.907.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.18%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'589.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").Value = "'144.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.96%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.503"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.32%  '
$ws.Range("D9").Value = "'2.906.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.10%  '
$ws.Range("D10").Value = "'6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.37%  '
$ws.Range("D11").Value = "'0.144"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.62%  '
$ws.Range("D12").Value = "'0.444"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.39%  '
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.72%  '
$ws.Range("D14").Value = "'33.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.04%  '
$ws.Range("D16").Value = "'3.387.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.21%  '
$ws.Range("D17").Value = "'60.744.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.88%  '
$ws.Range("D18").Value = "'6.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.58%  '
$ws.Range("D19").Value = "'2.906.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.27%  '
$ws.Range("D20").Value = "'429.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.02%  '
$ws.Range("D21").Value = "'13.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.26%  '
$ws.Range("D22").Value = "'0.683"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("D23").Value = "'7.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.15%  '
$ws.Range("D24").Value = "'81.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.85%  '
$ws.Range("D25").Value = "'10.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.13%  '
$ws.Range("D26").Value = "'2.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.54%  '
$ws.Range("D27").Value = "'11.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.94%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = "'2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D32").Value = "'7.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.08%  '
$ws.Range("D33").Value = "'26.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.45%  '
$ws.Range("E34").Value = '  -3.67%  '
$ws.Range("D35").Value = "'0.0₃0848"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.53%  '
$ws.Range("E36").Value = '  -3.28%  '
$ws.Range("E37").Value = '  -5.29%  '
$ws.Range("D38").Value = "'3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.37%  '
$ws.Range("D39").Value = "'49.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.40%  '
$ws.Range("E40").Value = '  -4.57%  '
$ws.Range("D41").Value = "'2.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.86%  '
$ws.Range("D42").Value = "'8.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.57%  '
$ws.Range("D43").Value = "'0.294"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.24%  '
$ws.Range("D44").Value = "'40.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.60%  '
$ws.Range("D45").Value = "'0.0350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.17%  '
$ws.Range("D46").Value = "'373.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.56%  '
$ws.Range("D47").Value = "'2.700.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("D48").Value = "'132.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").Value = "'24.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.12%  '
$ws.Range("E51").Value = '  -2.77%  '
